$d = $word.ActiveDocument

# 1. Update the status-line text inside the highlighted example run.
$d.Content.Find.Execute(
    "MW#2023/08/22 16:42#100%#Mounted#16GB#1GB#1%#8#1#",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "MW#2023/09/03 18:57#N/A#Mounted#16GB#1GB#6.25%#25.68#102.98#", 2)

# 2. The closing curly quote right after the status line now also needs
#    to be highlighted yellow, matching the preceding run.
$full = $d.Content.Text
$newStatus = "MW#2023/09/03 18:57#N/A#Mounted#16GB#1GB#6.25%#25.68#102.98#"
$statusIdx = $full.IndexOf($newStatus)
$quoteStart = $statusIdx + $newStatus.Length
$quoteRng = $d.Range($quoteStart, $quoteStart + 1)
$quoteRng.Font.HighlightColorIndex = 7
